$d = $word.ActiveDocument

# The document originally splits several title/author/abstract phrases
# into one run per word (with separate space-only runs in between).
# Collapse each of those paragraphs back down into a single run holding
# the full phrase, leaving every other paragraph untouched.

function Merge-ParagraphText($StyleName, $FullText) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Style.NameLocal -eq $StyleName) {
            $rng = $p.Range
            # Scope the Find/Replace to just this paragraph's range so we
            # never touch matching text elsewhere in the document.
            $rng.Find.Execute($FullText, $true, $false, $false, $false, $false,
                               $true, 1, $false, $FullText, 2)
            break
        }
    }
}

Merge-ParagraphText "Title" "Factsheet: Laws of logarithms"
Merge-ParagraphText "Author" "Millie Pike"
Merge-ParagraphText "Abstract" "A list of laws of logarithms."
